$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source report is appended daily; this update carries the series
# forward from 2021-09-21 through 2021-12-08 (commit: "aggiornamento fino
# a 8/12"), adding rows 386-464 with columns: date serial, nuovi pos.,
# somma mobile 7gg., somma mobile 7gg. per 100mila abitanti.

# Seed the new rows with the same look as the existing data (row 385):
# column A keeps the boxed/centered date style + number format, while
# B-D stay on the default style, matching every prior row in the sheet.
$ws.Range("A385:D385").Copy()
$ws.Range("A386:D464").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newData = @(
    @(44460, 0, 0, 0),
    @(44461, 0, 0, 0),
    @(44462, 0, 0, 0),
    @(44463, 0, 0, 0),
    @(44464, 0, 0, 0),
    @(44465, 0, 0, 0),
    @(44466, 0, 0, 0),
    @(44467, 0, 0, 0),
    @(44468, 0, 0, 0),
    @(44469, 0, 0, 0),
    @(44470, 0, 0, 0),
    @(44471, 0, 0, 0),
    @(44472, 0, 0, 0),
    @(44473, 0, 0, 0),
    @(44474, 0, 0, 0),
    @(44475, 0, 0, 0),
    @(44476, 0, 0, 0),
    @(44477, 0, 0, 0),
    @(44478, 0, 0, 0),
    @(44479, 0, 0, 0),
    @(44480, 0, 0, 0),
    @(44481, 0, 0, 0),
    @(44482, 0, 0, 0),
    @(44483, 0, 0, 0),
    @(44484, 0, 0, 0),
    @(44485, 0, 0, 0),
    @(44486, 0, 0, 0),
    @(44487, 0, 0, 0),
    @(44488, 0, 0, 0),
    @(44489, 0, 0, 0),
    @(44490, 0, 0, 0),
    @(44491, 0, 0, 0),
    @(44492, 0, 0, 0),
    @(44493, 1, 1, 33.71544167228591),
    @(44494, 0, 1, 33.71544167228591),
    @(44495, 1, 2, 67.43088334457181),
    @(44496, 0, 2, 67.43088334457181),
    @(44497, 1, 3, 101.1463250168577),
    @(44498, 1, 4, 134.8617666891436),
    @(44499, 1, 5, 168.5772083614295),
    @(44500, 0, 4, 134.8617666891436),
    @(44501, 0, 4, 134.8617666891436),
    @(44502, 1, 4, 134.8617666891436),
    @(44503, 0, 4, 134.8617666891436),
    @(44504, 0, 3, 101.1463250168577),
    @(44505, 1, 3, 101.1463250168577),
    @(44506, 1, 3, 101.1463250168577),
    @(44507, 0, 3, 101.1463250168577),
    @(44508, 0, 3, 101.1463250168577),
    @(44509, 2, 4, 134.8617666891436),
    @(44510, 0, 4, 134.8617666891436),
    @(44511, 1, 5, 168.5772083614295),
    @(44512, 2, 6, 202.2926500337155),
    @(44513, 0, 5, 168.5772083614295),
    @(44514, 1, 6, 202.2926500337155),
    @(44515, 5, 11, 370.8698583951449),
    @(44516, 0, 9, 303.4389750505732),
    @(44517, 0, 9, 303.4389750505732),
    @(44518, 4, 12, 404.5853000674309),
    @(44519, 3, 13, 438.3007417397168),
    @(44520, 2, 15, 505.7316250842886),
    @(44521, 0, 14, 472.0161834120027),
    @(44522, 3, 12, 404.5853000674309),
    @(44523, 0, 12, 404.5853000674309),
    @(44524, 2, 14, 472.0161834120027),
    @(44525, 3, 13, 438.3007417397168),
    @(44526, 0, 10, 337.1544167228591),
    @(44527, 3, 11, 370.8698583951449),
    @(44528, 2, 13, 438.3007417397168),
    @(44529, 2, 12, 404.5853000674309),
    @(44530, 0, 12, 404.5853000674309),
    @(44531, 2, 12, 404.5853000674309),
    @(44532, 6, 15, 505.7316250842886),
    @(44533, 3, 18, 606.8779501011463),
    @(44534, 1, 16, 539.4470667565745),
    @(44535, 0, 14, 472.0161834120027),
    @(44536, 1, 13, 438.3007417397168),
    @(44537, 4, 17, 573.1625084288604),
    @(44538, 0, 15, 505.7316250842886)
)

$startRow = 386
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $newData[$i]
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
}

